$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.356.58"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "1.839.01"
$ws.Range("E3").Value = "  -0.30%  "

$c = $ws.Range("D4")
$c.Value = "'0.9995"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.20%  "

$c = $ws.Range("D5")
$c.Value = "'238.69"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "

$c = $ws.Range("D6")
$c.Value = "'0.6263"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("E7").Value = "  +0.11%  "

$c = $ws.Range("D8")
$c.Value = "'0.07401"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.10%  "

$c = $ws.Range("D9")
$c.Value = "'0.2888"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "

$c = $ws.Range("D10")
$c.Value = "'24.90"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "

$c = $ws.Range("D11")
$c.Value = "'0.07716"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").Value = "1.839.04"
$ws.Range("E12").Value = "  -0.28%  "

$c = $ws.Range("D13")
$c.Value = "'4.967"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.37%  "

$c = $ws.Range("D14")
$c.Value = "'0.6734"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.93%  "

$c = $ws.Range("D15")
$c.Value = "'0.00001023"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.01%  "

$c = $ws.Range("D16")
$c.Value = "'81.66"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.35%  "

$c = $ws.Range("D17")
$c.Value = "'6.198"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("D18").Value = "29.415.88"
$ws.Range("E18").Value = "  +0.16%  "

$c = $ws.Range("D19")
$c.Value = "'232.64"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.64%  "

$c = $ws.Range("D20")
$c.Value = "'12.30"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("E21").Value = "  +0.16%  "

$c = $ws.Range("D22")
$c.Value = "'7.286"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.90%  "

$ws.Range("E23").Value = "  +0.25%  "

$c = $ws.Range("D24")
$c.Value = "'157.92"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.43%  "

$c = $ws.Range("D25")
$c.Value = "'8.485"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "

$c = $ws.Range("D26")
$c.Value = "'0.1343"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.74%  "

$c = $ws.Range("D27")
$c.Value = "'17.29"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.23%  "

$c = $ws.Range("D28")
$c.Value = "'0.07217"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +9.68%  "

$c = $ws.Range("D29")
$c.Value = "'1.469"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.73%  "

$c = $ws.Range("D30")
$c.Value = "'1.477"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.42%  "

$c = $ws.Range("D31")
$c.Value = "'4.033"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.73%  "

$c = $ws.Range("D32")
$c.Value = "'4.024"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.59%  "

$c = $ws.Range("D33")
$c.Value = "'1.814"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.59%  "

$c = $ws.Range("D34")
$c.Value = "'1.139"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.14%  "

$c = $ws.Range("D35")
$c.Value = "'0.6959"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "

$c = $ws.Range("D36")
$c.Value = "'2.573"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.33%  "

$c = $ws.Range("D37")
$c.Value = "'0.01839"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "

$c = $ws.Range("D38")
$c.Value = "'6.903"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.67%  "

$c = $ws.Range("D39")
$c.Value = "'2.815"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.95%  "

$ws.Range("D40").Value = "1.233.00"
$ws.Range("E40").Value = "  -2.40%  "

$c = $ws.Range("D41")
$c.Value = "'0.9575"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +4.39%  "

$c = $ws.Range("D42")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").Value = "2.015.32"
$ws.Range("E43").Value = "  +0.41%  "

$c = $ws.Range("D44")
$c.Value = "'100.68"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "

$c = $ws.Range("D45")
$c.Value = "'65.26"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.30%  "

$c = $ws.Range("D46")
$c.Value = "'0.00000000118"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.03%  "

$c = $ws.Range("D47")
$c.Value = "'1.709"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "

$c = $ws.Range("D48")
$c.Value = "'6.931"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.01%  "

$c = $ws.Range("D49")
$c.Value = "'8.828"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.43%  "

$c = $ws.Range("D50")
$c.Value = "'0.3892"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.40%  "

$c = $ws.Range("D51")
$c.Value = "'0.1129"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.70%  "

